# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 4 of the
# zh-cn and de-de report sheets, regenerating the handback report times.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-24 16:54:54"
$zhcn.Range("H4").Value = "2016-03-24 16:55:32"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-24 16:54:58"
$dede.Range("H4").Value = "2016-03-24 16:55:39"
